# Generated by Katalon AI
# Expands the "AI Generated" sheet from a single input_Name column into a
# six-column locator table (button_closeActions_class,
# button_closeActions_class_1, div_menuBackdrops_class,
# div_menuBackdrops_class_1, input_Name, p_sessionDetails_class) with the
# associated locator values underneath.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (A:F) -----------------------------------------------
# ColumnWidth in the object model is offset by ~0.8333 chars vs. the raw
# <col width="..."> stored in the sheet XML, so back that offset out to
# land on the exact target widths (128, 78, 53, 53, 12, 96).
$ws.Columns.Item(1).ColumnWidth = 127.16666666666667
$ws.Columns.Item(2).ColumnWidth = 77.16666666666667
$ws.Columns.Item(3).ColumnWidth = 52.166666666666664
$ws.Columns.Item(4).ColumnWidth = 52.166666666666664
$ws.Columns.Item(5).ColumnWidth = 11.166666666666666
$ws.Columns.Item(6).ColumnWidth = 95.16666666666667

# --- Row 1 : header labels ----------------------------------------------
$ws.Range("A1").Value = "button_closeActions_class"
$ws.Range("B1").Value = "button_closeActions_class_1"
$ws.Range("C1").Value = "div_menuBackdrops_class"
$ws.Range("D1").Value = "div_menuBackdrops_class_1"
$ws.Range("E1").Value = "input_Name"
$ws.Range("F1").Value = "p_sessionDetails_class"

# The original A1 carried the "Pandas" header style (bold / bordered /
# centered). Copy that formatting onto the new header cells so they all
# share the same style entry instead of Excel minting a new xf for each.
$ws.Range("A1").Copy()
$ws.Range("B1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 2 : locator values ----------------------------------------------
$ws.Range("A2").Value = 'go1475592160\ go1671063245"]:nth-child(1) [class="MuiButtonBase-root\ MuiIconButton-root\ MuiIconButton-sizeSmall\ css-1fubc2b'
$ws.Range("B2").Value = 'MuiButtonBase-root\ MuiIconButton-root\ MuiIconButton-sizeSmall\ css-1fubc2b'
$ws.Range("C2").Value = 'MuiBackdrop-root\ MuiBackdrop-invisible\ css-esi9ax'
$ws.Range("D2").Value = 'MuiBackdrop-root\ MuiBackdrop-invisible\ css-esi9ax'
$ws.Range("F2").Value = 'MuiBox-root\ css-0"]:nth-child(1) [class="MuiTypography-root\ MuiTypography-body1\ css-1lpm9pj'
# E2 (the old A2) stays blank, same as before the edit.
